$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.343.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "'3.323.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D5").Value = "'588.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "'183.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "'0.648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.13%  "
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'3.903.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("D14").Value = "'66.357.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "'3.305.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "'426.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").Value = "'5.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").Value = "'71.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'3.463.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  +7.39%  "
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'22.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("D35").Value = "'6.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("D37").Value = "'159.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").Value = "'2.889.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").Value = "'26.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").Value = "'40.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "'5.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").Value = "'23.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.78%  "
$ws.Range("D49").Value = "'314.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("E51").Value = "  +5.05%  "
